$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.366139666666667
$ws.Range("H2").Value = 10.098419
$ws.Range("I2").Value = 0.01725116351498256
$ws.Range("J2").Value = 0.01815407111703398
$ws.Range("M2").Value = 121.928739
$ws.Range("N2").Value = 365.786217
$ws.Range("O2").Value = 0.2282232151508951
$ws.Range("P2").Value = 0.2419720431319445
$ws.Range("Q2").Value = 410.429164854547
$ws.Range("R2").Value = 3693.862483690923
$ws.Range("S2").Value = 0.003937116002483136
$ws.Range("T2").Value = 0.004392777679351334
$ws.Range("G3").Value = 3.366139666666667
$ws.Range("H3").Value = 10.098419
$ws.Range("I3").Value = 0.01725116351498256
$ws.Range("J3").Value = 0.01815407111703398
$ws.Range("O3").Value = 0.2768624053389947
$ws.Range("P3").Value = 0.2935413991166814
$ws.Range("Q3").Value = 497.9002934814233
$ws.Range("R3").Value = 4481.10264133281
$ws.Range("S3").Value = 0.004776198625654378
$ws.Range("T3").Value = 0.00532897143535789
$ws.Range("G4").Value = 3.366139666666667
$ws.Range("H4").Value = 10.098419
$ws.Range("I4").Value = 0.01725116351498256
$ws.Range("J4").Value = 0.01815407111703398
$ws.Range("M4").Value = 83.50496933333334
$ws.Range("N4").Value = 250.514908
$ws.Range("O4").Value = 0.1563025480180701
$ws.Range("P4").Value = 0.1657186665504434
$ws.Range("Q4").Value = 281.0893896367169
$ws.Range("R4").Value = 2529.804506730452
$ws.Range("S4").Value = 0.00269640081366814
$ws.Range("T4").Value = 0.003008468457976789
$ws.Range("G5").Value = 3.366139666666667
$ws.Range("H5").Value = 10.098419
$ws.Range("I5").Value = 0.01725116351498256
$ws.Range("J5").Value = 0.01815407111703398
$ws.Range("M5").Value = 91.06846250000001
$ws.Range("N5").Value = 182.136925
$ws.Range("O5").Value = 0.1704597085236707
$ws.Range("P5").Value = 0.1204857969594293
$ws.Range("Q5").Value = 306.5491640035959
$ws.Range("R5").Value = 1839.294984021575
$ws.Range("S5").Value = 0.00294062830445811
$ws.Range("T5").Value = 0.002187307726593996
$ws.Range("G6").Value = 3.366139666666667
$ws.Range("H6").Value = 10.098419
$ws.Range("I6").Value = 0.01725116351498256
$ws.Range("J6").Value = 0.01815407111703398
$ws.Range("M6").Value = 89.83562999999999
$ws.Range("N6").Value = 269.50689
$ws.Range("O6").Value = 0.1681521229683693
$ws.Range("P6").Value = 0.1782820942415013
$ws.Range("Q6").Value = 302.39927762299
$ws.Range("R6").Value = 2721.59349860691
$ws.Range("S6").Value = 0.002900819768718793
$ws.Range("T6").Value = 0.003236545817753968
$ws.Range("I7").Value = 0.7504462978934635
$ws.Range("J7").Value = 0.7897238612132288
$ws.Range("M7").Value = 121.928739
$ws.Range("N7").Value = 365.786217
$ws.Range("O7").Value = 0.2282232151508951
$ws.Range("P7").Value = 0.2419720431319445
$ws.Range("Q7").Value = 17854.16079588485
$ws.Range("R7").Value = 160687.4471629637
$ws.Range("S7").Value = 0.1712692669033326
$ws.Range("T7").Value = 0.1910910962078132
$ws.Range("I8").Value = 0.7504462978934635
$ws.Range("J8").Value = 0.7897238612132288
$ws.Range("O8").Value = 0.2768624053389947
$ws.Range("P8").Value = 0.2935413991166814
$ws.Range("Q8").Value = 21659.25977333018
$ws.Range("S8").Value = 0.2077703671125281
$ws.Range("T8").Value = 0.2318166471363591
$ws.Range("I9").Value = 0.7504462978934635
$ws.Range("J9").Value = 0.7897238612132288
$ws.Range("M9").Value = 83.50496933333334
$ws.Range("N9").Value = 250.514908
$ws.Range("O9").Value = 0.1563025480180701
$ws.Range("P9").Value = 0.1657186665504434
$ws.Range("Q9").Value = 12227.72548917091
$ws.Range("R9").Value = 110049.5294025382
$ws.Range("S9").Value = 0.117296668511476
$ws.Range("T9").Value = 0.1308719852233237
$ws.Range("I10").Value = 0.7504462978934635
$ws.Range("J10").Value = 0.7897238612132288
$ws.Range("M10").Value = 91.06846250000001
$ws.Range("N10").Value = 182.136925
$ws.Range("O10").Value = 0.1704597085236707
$ws.Range("P10").Value = 0.1204857969594293
$ws.Range("Q10").Value = 13335.25620164915
$ws.Range("R10").Value = 80011.5372098949
$ws.Range("S10").Value = 0.1279208572015876
$ws.Range("T10").Value = 0.09515050879615364
$ws.Range("I11").Value = 0.7504462978934635
$ws.Range("J11").Value = 0.7897238612132288
$ws.Range("M11").Value = 89.83562999999999
$ws.Range("N11").Value = 269.50689
$ws.Range("O11").Value = 0.1681521229683693
$ws.Range("P11").Value = 0.1782820942415013
$ws.Range("Q11").Value = 13154.73116817535
$ws.Range("R11").Value = 118392.5805135782
$ws.Range("S11").Value = 0.1261891381645392
$ws.Range("T11").Value = 0.1407936238495791
$ws.Range("G12").Value = 14.89209833333333
$ws.Range("H12").Value = 44.676295
$ws.Range("I12").Value = 0.07632066665966204
$ws.Range("J12").Value = 0.08031520940808551
$ws.Range("M12").Value = 121.928739
$ws.Range("N12").Value = 365.786217
$ws.Range("O12").Value = 0.2282232151508951
$ws.Range("P12").Value = 0.2419720431319445
$ws.Range("Q12").Value = 1815.774770847335
$ws.Range("R12").Value = 16341.97293762601
$ws.Range("S12").Value = 0.0174181479275278
$ws.Range("T12").Value = 0.01943403531504442
$ws.Range("G13").Value = 14.89209833333333
$ws.Range("H13").Value = 44.676295
$ws.Range("I13").Value = 0.07632066665966204
$ws.Range("J13").Value = 0.08031520940808551
$ws.Range("O13").Value = 0.2768624053389947
$ws.Range("P13").Value = 0.2935413991166814
$ws.Range("Q13").Value = 2202.754747269116
$ws.Range("R13").Value = 19824.79272542204
$ws.Range("S13").Value = 0.02113032334846965
$ws.Range("T13").Value = 0.02357583893999867
$ws.Range("G14").Value = 14.89209833333333
$ws.Range("H14").Value = 44.676295
$ws.Range("I14").Value = 0.07632066665966204
$ws.Range("J14").Value = 0.08031520940808551
$ws.Range("M14").Value = 83.50496933333334
$ws.Range("N14").Value = 250.514908
$ws.Range("O14").Value = 0.1563025480180701
$ws.Range("P14").Value = 0.1657186665504434
$ws.Range("Q14").Value = 1243.564214633984
$ws.Range("R14").Value = 11192.07793170586
$ws.Range("S14").Value = 0.01192911466534295
$ws.Range("T14").Value = 0.01330972940682756
$ws.Range("G15").Value = 14.89209833333333
$ws.Range("H15").Value = 44.676295
$ws.Range("I15").Value = 0.07632066665966204
$ws.Range("J15").Value = 0.08031520940808551
$ws.Range("M15").Value = 91.06846250000001
$ws.Range("N15").Value = 182.136925
$ws.Range("O15").Value = 0.1704597085236707
$ws.Range("P15").Value = 0.1204857969594293
$ws.Range("Q15").Value = 1356.200498615479
$ws.Range("R15").Value = 8137.202991692875
$ws.Range("S15").Value = 0.01300959859313823
$ws.Range("T15").Value = 0.009676842013496639
$ws.Range("G16").Value = 14.89209833333333
$ws.Range("H16").Value = 44.676295
$ws.Range("I16").Value = 0.07632066665966204
$ws.Range("J16").Value = 0.08031520940808551
$ws.Range("M16").Value = 89.83562999999999
$ws.Range("N16").Value = 269.50689
$ws.Range("O16").Value = 0.1681521229683693
$ws.Range("P16").Value = 0.1782820942415013
$ws.Range("Q16").Value = 1337.84103579695
$ws.Range("R16").Value = 12040.56932217255
$ws.Range("S16").Value = 0.01283348212518341
$ws.Range("T16").Value = 0.01431876373271821
$ws.Range("G17").Value = 29.1141605
$ws.Range("H17").Value = 58.22832099999999
$ws.Range("I17").Value = 0.1492074581338761
$ws.Range("J17").Value = 0.1046778788302885
$ws.Range("M17").Value = 121.928739
$ws.Range("N17").Value = 365.786217
$ws.Range("O17").Value = 0.2282232151508951
$ws.Range("P17").Value = 0.2419720431319445
$ws.Range("Q17").Value = 3549.852876808609
$ws.Range("R17").Value = 21299.11726085165
$ws.Range("S17").Value = 0.03405260581980578
$ws.Range("T17").Value = 0.02532912021128303
$ws.Range("G18").Value = 29.1141605
$ws.Range("H18").Value = 58.22832099999999
$ws.Range("I18").Value = 0.1492074581338761
$ws.Range("J18").Value = 0.1046778788302885
$ws.Range("O18").Value = 0.2768624053389947
$ws.Range("P18").Value = 0.2935413991166814
$ws.Range("Q18").Value = 4306.401543869964
$ws.Range("R18").Value = 25838.40926321979
$ws.Range("S18").Value = 0.04130993575346229
$ws.Range("T18").Value = 0.03072729100840933
$ws.Range("G19").Value = 29.1141605
$ws.Range("H19").Value = 58.22832099999999
$ws.Range("I19").Value = 0.1492074581338761
$ws.Range("J19").Value = 0.1046778788302885
$ws.Range("M19").Value = 83.50496933333334
$ws.Range("N19").Value = 250.514908
$ws.Range("O19").Value = 0.1563025480180701
$ws.Range("P19").Value = 0.1657186665504434
$ws.Range("Q19").Value = 2431.177079718244
$ws.Range("R19").Value = 14587.06247830947
$ws.Range("S19").Value = 0.02332150588962435
$ws.Range("T19").Value = 0.01734707849708429
$ws.Range("G20").Value = 29.1141605
$ws.Range("H20").Value = 58.22832099999999
$ws.Range("I20").Value = 0.1492074581338761
$ws.Range("J20").Value = 0.1046778788302885
$ws.Range("M20").Value = 91.06846250000001
$ws.Range("N20").Value = 182.136925
$ws.Range("O20").Value = 0.1704597085236707
$ws.Range("P20").Value = 0.1204857969594293
$ws.Range("Q20").Value = 2651.381833713231
$ws.Range("R20").Value = 10605.52733485292
$ws.Range("S20").Value = 0.02543385982305833
$ws.Range("T20").Value = 0.01261219765488988
$ws.Range("G21").Value = 29.1141605
$ws.Range("H21").Value = 58.22832099999999
$ws.Range("I21").Value = 0.1492074581338761
$ws.Range("J21").Value = 0.1046778788302885
$ws.Range("M21").Value = 89.83562999999999
$ws.Range("N21").Value = 269.50689
$ws.Range("O21").Value = 0.1681521229683693
$ws.Range("P21").Value = 0.1782820942415013
$ws.Range("Q21").Value = 2615.488950438615
$ws.Range("R21").Value = 15692.93370263169
$ws.Range("S21").Value = 0.02508955084792535
$ws.Range("T21").Value = 0.01866219145862194
$ws.Range("G22").Value = 1.32186
$ws.Range("H22").Value = 3.96558
$ws.Range("I22").Value = 0.006774413798015763
$ws.Range("J22").Value = 0.007128979431363227
$ws.Range("M22").Value = 121.928739
$ws.Range("N22").Value = 365.786217
$ws.Range("O22").Value = 0.2282232151508951
$ws.Range("P22").Value = 0.2419720431319445
$ws.Range("Q22").Value = 161.17272293454
$ws.Range("R22").Value = 1450.55450641086
$ws.Range("S22").Value = 0.001546078497745744
$ws.Range("T22").Value = 0.001725013718452568
$ws.Range("G23").Value = 1.32186
$ws.Range("H23").Value = 3.96558
$ws.Range("I23").Value = 0.006774413798015763
$ws.Range("J23").Value = 0.007128979431363227
$ws.Range("O23").Value = 0.2768624053389947
$ws.Range("P23").Value = 0.2935413991166814
$ws.Range("Q23").Value = 195.5220362538
$ws.Range("R23").Value = 1759.6983262842
$ws.Range("S23").Value = 0.001875580498880319
$ws.Range("T23").Value = 0.002092650596556406
$ws.Range("G24").Value = 1.32186
$ws.Range("H24").Value = 3.96558
$ws.Range("I24").Value = 0.006774413798015763
$ws.Range("J24").Value = 0.007128979431363227
$ws.Range("M24").Value = 83.50496933333334
$ws.Range("N24").Value = 250.514908
$ws.Range("O24").Value = 0.1563025480180701
$ws.Range("P24").Value = 0.1657186665504434
$ws.Range("Q24").Value = 110.38187876296
$ws.Range("R24").Value = 993.43690886664
$ws.Range("S24").Value = 0.001058858137958635
$ws.Range("T24").Value = 0.001181404965231052
$ws.Range("G25").Value = 1.32186
$ws.Range("H25").Value = 3.96558
$ws.Range("I25").Value = 0.006774413798015763
$ws.Range("J25").Value = 0.007128979431363227
$ws.Range("M25").Value = 91.06846250000001
$ws.Range("N25").Value = 182.136925
$ws.Range("O25").Value = 0.1704597085236707
$ws.Range("P25").Value = 0.1204857969594293
$ws.Range("Q25").Value = 120.37975784025
$ws.Range("R25").Value = 722.2785470415001
$ws.Range("S25").Value = 0.0011547646014285
$ws.Range("T25").Value = 0.0008589407682951778
$ws.Range("G26").Value = 1.32186
$ws.Range("H26").Value = 3.96558
$ws.Range("I26").Value = 0.006774413798015763
$ws.Range("J26").Value = 0.007128979431363227
$ws.Range("M26").Value = 89.83562999999999
$ws.Range("N26").Value = 269.50689
$ws.Range("O26").Value = 0.1681521229683693
$ws.Range("P26").Value = 0.1782820942415013
$ws.Range("Q26").Value = 118.7501258718
$ws.Range("R26").Value = 1068.7511328462
$ws.Range("S26").Value = 0.001139132062002564
$ws.Range("T26").Value = 0.001270969382828023
